$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-05 23:17:47"
$ws.Range("H2").Value = "96%"
$ws.Range("E3").Value = "2026-02-05 23:17:50"
$ws.Range("E4").Value = "2026-02-05 23:17:52"
$ws.Range("J4").Value = "989.5 hPa"
$ws.Range("O4").Value = "11.9 °C"
$ws.Range("E5").Value = "2026-02-05 23:17:54"
$ws.Range("E6").Value = "2026-02-05 23:17:57"
$ws.Range("O6").Value = "13.3 °C"
$ws.Range("E7").Value = "2026-02-05 23:18:00"
$ws.Range("E8").Value = "2026-02-05 23:18:02"
$ws.Range("E9").Value = "2026-02-05 23:18:04"
$ws.Range("E10").Value = "2026-02-05 23:18:07"
$ws.Range("E11").Value = "2026-02-05 23:18:09"
$ws.Range("E12").Value = "2026-02-05 23:18:12"
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-05 23:18:14"
$ws.Range("E14").Value = "2026-02-05 23:18:16"
$ws.Range("E15").Value = "2026-02-05 23:18:19"
$ws.Range("K15").Value = "6.6 MJ/m2"
$ws.Range("O15").Value = "9.2 °C"
$ws.Range("E16").Value = "2026-02-05 23:18:21"
$ws.Range("E17").Value = "2026-02-05 23:18:24"
$ws.Range("H17").Value = "98%"
$ws.Range("I17").Value = "9.1 mm"
$ws.Range("M17").Value = "3.4 °C 22:38 TU"
$ws.Range("E18").Value = "2026-02-05 23:18:26"
$ws.Range("I18").Value = "2.8 mm"
$ws.Range("E19").Value = "2026-02-05 23:18:29"
$ws.Range("J19").Value = "992.7 hPa"
$ws.Range("E20").Value = "2026-02-05 23:18:31"
$ws.Range("E21").Value = "2026-02-05 23:18:34"
$ws.Range("E22").Value = "2026-02-05 23:18:36"
$ws.Range("E23").Value = "2026-02-05 23:18:39"
$ws.Range("E24").Value = "2026-02-05 23:18:41"
$ws.Range("O24").Value = "10.7 °C"
$ws.Range("E25").Value = "2026-02-05 23:18:43"
$ws.Range("J25").Value = "993.8 hPa"
$ws.Range("O25").Value = "0.9 °C"
$ws.Range("E26").Value = "2026-02-05 23:18:46"
$ws.Range("E27").Value = "2026-02-05 23:18:49"
$ws.Range("O27").Value = "8.7 °C"
$ws.Range("E28").Value = "2026-02-05 23:18:51"
$ws.Range("O28").Value = "3.1 °C"
$ws.Range("E29").Value = "2026-02-05 23:18:53"
$ws.Range("H29").Value = "76%"
$ws.Range("O29").Value = "9.9 °C"
$ws.Range("E30").Value = "2026-02-05 23:18:56"
$ws.Range("E31").Value = "2026-02-05 23:18:58"
$ws.Range("I31").Value = "20.5 mm"
$ws.Range("M31").Value = "5.2 °C 22:42 TU"
$ws.Range("O31").Value = "3.8 °C"
$ws.Range("E32").Value = "2026-02-05 23:19:01"
$ws.Range("H32").Value = "78%"
$ws.Range("E33").Value = "2026-02-05 23:19:03"
$ws.Range("H33").Value = "83%"
$ws.Range("E34").Value = "2026-02-05 23:19:06"
$ws.Range("O34").Value = "4.8 °C"
$ws.Range("E35").Value = "2026-02-05 23:19:08"
$ws.Range("E36").Value = "2026-02-05 23:19:11"
